$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 3.09825757489699
$ws.Range("E2").Value = 7.865470614547343

$ws.Range("C3").Value = -1.791203563722299
$ws.Range("E3").Value = -6.760862998203643

$ws.Range("C4").Value = 0.5799958470386946
$ws.Range("E4").Value = 5.643342995751777

$ws.Range("C5").Value = 3.181454202131073
$ws.Range("E5").Value = 4.38978860149748

$ws.Range("C6").Value = 0.5930547804883668
$ws.Range("E6").Value = -1.194610791900008

$ws.Range("C7").Value = -0.3951783438669754
$ws.Range("E7").Value = 0.03694906323863378

$ws.Range("C8").Value = 3.292216014290039
$ws.Range("E8").Value = 7.617133650412211

$ws.Range("C9").Value = 1.670328650030184
$ws.Range("E9").Value = 2.037906845818616

$ws.Range("C10").Value = 2.562791874943371
$ws.Range("E10").Value = 3.265947405805814

$ws.Range("C11").Value = 1.526411006965533
$ws.Range("E11").Value = 0.6601843988560674

$ws.Range("C12").Value = 1.63465618619294
$ws.Range("E12").Value = 1.551857746372698

$ws.Range("C13").Value = 1.35261353265177
$ws.Range("E13").Value = 0.8024032016000104

$ws.Range("C14").Value = -2.082763426755907
$ws.Range("E14").Value = -5.866344937500023

$ws.Range("C15").Value = -0.1380317107957718
$ws.Range("E15").Value = 7.749494937649115

$ws.Range("C16").Value = 3.848999231984762
$ws.Range("E16").Value = 2.866003071127765

$ws.Range("C17").Value = -0.3745803349312071
$ws.Range("E17").Value = 0.645722451525943

$ws.Range("C18").Value = -1.091476630333243
$ws.Range("E18").Value = -0.4907904687545206

$ws.Range("C19").Value = 1.626992717807862
$ws.Range("E19").Value = 0.6270138473519316
